# "remove GTIN because of improper format"
# The K column (GTIN) on the "Items" sheet was stored as numeric values in
# cells formatted as text ("@") - an improper format for a GTIN/barcode
# value. Clear those stored values, keeping the existing text formatting.

$wb = $excel.ActiveWorkbook

$items = $wb.Worksheets.Item("Items")
$items.Range("K7:K107").ClearContents()

# Leave the selection / cursor roughly where the author left it afterwards
# (the "Items" sheet stays the active sheet/tab).
$items.Activate()
$items.Range("K7:K119").Select()
